{"js": "// Update the worksheet date and all 25 three-digit-by-one-digit\n// multiplication prompts to the next day's generated values.\nconst replacements = [\n  [\"2024-09-11 Wednesday\", \"2024-09-12 Thursday\"],\n  [\"681\u00d77=\", \"207\u00d78=\"],\n  [\"753\u00d78=\", \"489\u00d75=\"],\n  [\"320\u00d74=\", \"369\u00d78=\"],\n  [\"704\u00d75=\", \"487\u00d72=\"],\n  [\"443\u00d79=\", \"584\u00d75=\"],\n  [\"938\u00d75=\", \"177\u00d78=\"],\n  [\"845\u00d77=\", \"309\u00d79=\"],\n  [\"998\u00d72=\", \"607\u00d76=\"],\n  [\"381\u00d78=\", \"625\u00d75=\"],\n  [\"624\u00d76=\", \"636\u00d73=\"],\n  [\"637\u00d73=\", \"656\u00d76=\"],\n  [\"128\u00d77=\", \"305\u00d74=\"],\n  [\"551\u00d76=\", \"266\u00d77=\"],\n  [\"610\u00d73=\", \"616\u00d74=\"],\n  [\"845\u00d75=\", \"219\u00d72=\"],\n  [\"379\u00d76=\", \"785\u00d78=\"],\n  [\"937\u00d79=\", \"117\u00d75=\"],\n  [\"961\u00d72=\", \"266\u00d75=\"],\n  [\"901\u00d75=\", \"838\u00d78=\"],\n  [\"736\u00d78=\", \"306\u00d76=\"],\n  [\"846\u00d73=\", \"968\u00d75=\"],\n  [\"973\u00d75=\", \"388\u00d77=\"],\n  [\"997\u00d77=\", \"824\u00d72=\"],\n  [\"909\u00d72=\", \"114\u00d73=\"],\n  [\"296\u00d77=\", \"622\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 three-digit-by-one-digit\n# multiplication prompts to the next day's generated values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-09-11 Wednesday\", \"2024-09-12 Thursday\"),\n    @(\"681\u00d77=\", \"207\u00d78=\"),\n    @(\"753\u00d78=\", \"489\u00d75=\"),\n    @(\"320\u00d74=\", \"369\u00d78=\"),\n    @(\"704\u00d75=\", \"487\u00d72=\"),\n    @(\"443\u00d79=\", \"584\u00d75=\"),\n    @(\"938\u00d75=\", \"177\u00d78=\"),\n    @(\"845\u00d77=\", \"309\u00d79=\"),\n    @(\"998\u00d72=\", \"607\u00d76=\"),\n    @(\"381\u00d78=\", \"625\u00d75=\"),\n    @(\"624\u00d76=\", \"636\u00d73=\"),\n    @(\"637\u00d73=\", \"656\u00d76=\"),\n    @(\"128\u00d77=\", \"305\u00d74=\"),\n    @(\"551\u00d76=\", \"266\u00d77=\"),\n    @(\"610\u00d73=\", \"616\u00d74=\"),\n    @(\"845\u00d75=\", \"219\u00d72=\"),\n    @(\"379\u00d76=\", \"785\u00d78=\"),\n    @(\"937\u00d79=\", \"117\u00d75=\"),\n    @(\"961\u00d72=\", \"266\u00d75=\"),\n    @(\"901\u00d75=\", \"838\u00d78=\"),\n    @(\"736\u00d78=\", \"306\u00d76=\"),\n    @(\"846\u00d73=\", \"968\u00d75=\"),\n    @(\"973\u00d75=\", \"388\u00d77=\"),\n    @(\"997\u00d77=\", \"824\u00d72=\"),\n    @(\"909\u00d72=\", \"114\u00d73=\"),\n    @(\"296\u00d77=\", \"622\u00d74=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
